$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Gomez)
$ws.Range("F2").Value = 6.224379325864566
$ws.Range("G2").Value = 8.928571428571429

# Row 3 (Bravo)
$ws.Range("F3").Value = 6.224379325864566

# Row 6 (Arredondo) - Categoria change
$ws.Range("D6").Value = "Master-Mañanas"

# Row 10 (Cisternas)
$ws.Range("G10").Value = 10.71428571428572

# Row 13 (Boettiger)
$ws.Range("G13").Value = 80.35714285714286
